$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 6 of the delivery table. The shared-string table in the
# target workbook has its three brand-new strings appended in the order
# I6, B6, A6 (i.e. not simple left-to-right column order), so we set
# those cells in that exact order to reproduce the shared string table.
$ws.Range("I6").Value = "72 SCRIMSHIRE Lane"
$ws.Range("B6").Value = "QUITZON, Luettgen and Waters"
$ws.Range("A6").Value = "PO0024616-1"

$ws.Range("C6").Value = "Viola"
$ws.Range("D6").Value = "String"
$ws.Range("E6").Value = 2000
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 4000
$ws.Range("H6").Value = "GBP"
$ws.Range("J6").Value = "Edinburgh"
$ws.Range("K6").Value = "EH7 4GT"
$ws.Range("L6").Value = "United Kingdom"
$ws.Range("M6").Value = "+44 115 496 0157"
$ws.Range("N6").Value = "Credit"
$ws.Range("O6").Value = "PO0024697-20210127"

# P6 already carries the date number format (style index 1); just set
# the serial date value for 2021-01-27.
$ws.Range("P6").Value = 44223

# Column P (16th column) gets a best-fit custom width of 10.5 once it
# holds real dates.
$ws.Columns.Item(16).ColumnWidth = 9.67

# Selection moves to the first cell of the newly-populated row.
$ws.Range("A6").Select() | Out-Null
